# Update cryptos list (prices & 1h volume %) per upstream refresh.
# For Price (column D) values that look like plain numbers, a leading
# apostrophe forces Excel to store them as literal text (matching the
# source data, which mixes plain and dotted/thousands-style numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.257.00'
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").Value = '2.267.99'
$ws.Range("E3").Value = '  -0.54%  '

$ws.Range("D5").Value = '''307.31'
$ws.Range("E5").Value = '  +0.88%  '

$ws.Range("D6").Value = '''97.46'
$ws.Range("E6").Value = '  +3.75%  '

$ws.Range("D7").Value = '''0.526'
$ws.Range("E7").Value = '  -0.92%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '''0.497'
$ws.Range("E9").Value = '  +1.45%  '

$ws.Range("D10").Value = '''35.34'
$ws.Range("E10").Value = '  +3.99%  '

$ws.Range("D11").Value = '''0.0791'
$ws.Range("E11").Value = '  -1.63%  '

$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("E13").Value = '  +3.16%  '

$ws.Range("D14").Value = '2.619.13'
$ws.Range("E14").Value = '  -0.57%  '

$ws.Range("D15").Value = '''14.75'
$ws.Range("E15").Value = '  +2.82%  '

$ws.Range("D16").Value = '2.256.38'
$ws.Range("E16").Value = '  -0.93%  '

$ws.Range("D17").Value = '''0.796'
$ws.Range("E17").Value = '  +0.70%  '

$ws.Range("D18").Value = '42.096.83'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").Value = '''12.48'
$ws.Range("E19").Value = '  -2.62%  '

$ws.Range("D20").Value = '0.0₃0907'
$ws.Range("E20").Value = '  -1.12%  '

$ws.Range("D21").Value = '''6.03'
$ws.Range("E21").Value = '  +0.48%  '

$ws.Range("D22").Value = '''68.34'
$ws.Range("E22").Value = '  +0.39%  '

$ws.Range("D23").Value = '''238.52'
$ws.Range("E23").Value = '  -2.15%  '

$ws.Range("D24").Value = '''2.59'
$ws.Range("E24").Value = '  -0.24%  '

$ws.Range("D25").Value = '''1.95'
$ws.Range("E25").Value = '  -0.32%  '

$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("D27").Value = '''23.66'
$ws.Range("E27").Value = '  -1.79%  '

$ws.Range("D28").Value = '''37.62'
$ws.Range("E28").Value = '  +5.11%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '''9.51'
$ws.Range("E29").Value = '  -1.62%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '''2.12'
$ws.Range("E30").Value = '  +0.87%  '

$ws.Range("D31").Value = '''160.97'
$ws.Range("E31").Value = '  +0.45%  '

$ws.Range("D32").Value = '''5.24'
$ws.Range("E32").Value = '  -2.03%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("D34").Value = '''3.18'
$ws.Range("E34").Value = '  +3.10%  '

$ws.Range("D35").Value = '''0.0739'
$ws.Range("E35").Value = '  -1.79%  '

$ws.Range("D36").Value = '''17.30'
$ws.Range("E36").Value = '  +1.54%  '

$ws.Range("E37").Value = '  -0.37%  '

$ws.Range("D38").Value = '''0.105'
$ws.Range("E38").Value = '  -2.69%  '

$ws.Range("D39").Value = '''1.83'
$ws.Range("E39").Value = '  +0.78%  '

$ws.Range("E40").Value = '  -1.49%  '

$ws.Range("D41").Value = '''4.01'
$ws.Range("E41").Value = '  -4.58%  '

$ws.Range("D42").Value = '''2.32'
$ws.Range("E42").Value = '  +2.02%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''19.23'
$ws.Range("E43").Value = '  -2.91%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.943.61'
$ws.Range("E44").Value = '  -3.98%  '

$ws.Range("D45").Value = '''0.0282'
$ws.Range("E45").Value = '  -0.37%  '

$ws.Range("D46").Value = '''10.01'
$ws.Range("E46").Value = '  -2.26%  '

$ws.Range("D47").Value = '''2.89'
$ws.Range("E47").Value = '  -1.66%  '

$ws.Range("D48").Value = '''53.46'
$ws.Range("E48").Value = '  -0.18%  '

$ws.Range("D49").Value = '''71.90'
$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("D50").Value = '''92.12'
$ws.Range("E50").Value = '  +0.18%  '

$ws.Range("E51").Value = '  -1.53%  '
